$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 22:52"

# --- Reorder: Brasil now appears before Belgica (rows 14/15) ---
# Row 14 becomes Brasil with refreshed figures; row 15 becomes Belgica keeping its prior figures.
$ws.Range("A14").Value = "Brasil"
$ws.Range("B14").Value = 41325
$ws.Range("C14").Value = 582
$ws.Range("D14").Value = 22991
$ws.Range("E14").Value = 15705
$ws.Range("F14").Value = 8318
$ws.Range("G14").Value = 42
$ws.Range("H14").Value = 2629

$ws.Range("A15").Value = "Belgica"
$ws.Range("B15").Value = 40956
$ws.Range("C15").Value = 973
$ws.Range("D15").Value = 9002
$ws.Range("E15").Value = 25956
$ws.Range("F15").Value = 1079
$ws.Range("G15").Value = 170
$ws.Range("H15").Value = 5998

# Alemania (row 8) data refresh
$ws.Range("B8").Value = 148226
$ws.Range("C8").Value = 1161
$ws.Range("E8").Value = 48002

# --- Reorder: Costa de Marfil now appears before Tunez (rows 87/88) ---
$ws.Range("A87").Value = "Costa de Marfil"
$ws.Range("B87").Value = 916
$ws.Range("C87").Value = 37
$ws.Range("D87").Value = 303
$ws.Range("E87").Value = 600
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 13

$ws.Range("A88").Value = "Tunez"
$ws.Range("B88").Value = 884
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 148
$ws.Range("E88").Value = 698
$ws.Range("F88").Value = 35
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 38

# --- Reorder: Trinidad yTobago now appears before Etiopia (rows 139/140) ---
$ws.Range("A139").Value = "Trinidad yTobago"
$ws.Range("B139").Value = 115
$ws.Range("C139").Value = 1
$ws.Range("D139").Value = 28
$ws.Range("E139").Value = 79
$ws.Range("H139").Value = 8

$ws.Range("A140").Value = "Etiopia"
$ws.Range("B140").Value = 114
$ws.Range("C140").Value = 3
$ws.Range("D140").Value = 16
$ws.Range("E140").Value = 95
$ws.Range("H140").Value = 3

# Togo (row 147) data refresh
$ws.Range("B147").Value = 86
$ws.Range("C147").Value = 2
$ws.Range("D147").Value = 56
$ws.Range("E147").Value = 24

# Zimbabue (row 176) data refresh
$ws.Range("B176").Value = 28
$ws.Range("C176").Value = 3
$ws.Range("E176").Value = 23

# San Bartolome (row 208) data refresh
$ws.Range("D208").Value = 6
$ws.Range("E208").Value = 0
